$wb = $excel.ActiveWorkbook

# Duplicate the "basic" sheet (same A:D layout/formulas) to the end of the
# workbook, producing the new "tournament5" machine sheet.
$basic = $wb.Worksheets.Item("basic")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$basic.Copy($null, $lastSheet)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "tournament5"

# Replace the slot symbols (Simpsons family) and their reel/payout data.
$ws.Range("A2").Value = "Maggie"
$ws.Range("B2").Value = 6
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = 24

$ws.Range("A3").Value = "Lisa"
$ws.Range("B3").Value = 8
$ws.Range("C3").Value = 6
$ws.Range("D3").Value = 6

$ws.Range("A4").Value = "Marge"
$ws.Range("B4").Value = 8
$ws.Range("C4").Value = 6
$ws.Range("D4").Value = 6

$ws.Range("A5").Value = "Bart"
$ws.Range("B5").Value = 12
$ws.Range("C5").Value = 6
$ws.Range("D5").Value = 2

$ws.Range("A6").Value = "Homer"
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 4
$ws.Range("D6").Value = 1

# Winning combination labels (entered top-down for the single-symbol rows,
# then bottom-up for the Homer rows).
$ws.Range("A11").Value = "Maggie+Any+Any"
$ws.Range("A12").Value = "Maggie+Maggie+Any"
$ws.Range("A13").Value = "Lisa+Lisa+Lisa"
$ws.Range("A14").Value = "Marge+Marge+Marge"
$ws.Range("A15").Value = "Bart+Bart+Bart"
$ws.Range("A18").Value = "Homer+Homer+Homer"
$ws.Range("A17").Value = "Homer+Homer+Any"
$ws.Range("A16").Value = "Homer+Any+Any"

# Updated payout multipliers.
$ws.Range("D11").Value = 2
$ws.Range("D12").Value = 5
$ws.Range("D13").Value = 10
$ws.Range("D14").Value = 12
$ws.Range("D15").Value = 15
$ws.Range("D16").Value = 5
$ws.Range("D17").Value = 25
$ws.Range("D18").Value = 50

# Incidental cursor move on "tournament" left over from the editing
# session, then land back on the new sheet's first data cell.
$t1 = $wb.Worksheets.Item("tournament")
$t1.Activate()
$t1.Range("C15").Select()
$ws.Activate()
$ws.Range("A2").Select()
